$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 156, pushing the existing "1a (cosecha)" /
# "2a (cosecha)" rows down to 158/159, and fill the freed rows 156/157 with
# the new "1a nueva(o)" / "2a nueva(o)" price entries.
$ws.Rows.Item(156).Resize(2).Insert()

# Row 156 - "1a nueva(o)"
$ws.Cells.Item(156, 1).Value = 11
$ws.Cells.Item(156, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(156, 3).Value = "Bíobío"
$ws.Cells.Item(156, 4).Value = 44568
$ws.Cells.Item(156, 4).Style = $ws.Cells.Item(158, 4).Style
$ws.Cells.Item(156, 4).NumberFormat = $ws.Cells.Item(158, 4).NumberFormat
$ws.Cells.Item(156, 5).Value = 8
$ws.Cells.Item(156, 6).Value = 100112045
$ws.Cells.Item(156, 7).Value = "Zapallo"
$ws.Cells.Item(156, 8).Value = "Camote"
$ws.Cells.Item(156, 9).Value = "1a nueva(o)"
$ws.Cells.Item(156, 10).Value = 1000
$ws.Cells.Item(156, 11).Value = 450
$ws.Cells.Item(156, 12).Value = 500
$ws.Cells.Item(156, 13).Value = 475
$ws.Cells.Item(156, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(156, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(156, 16).Value = 475
$ws.Cells.Item(156, 17).Value = 1
$ws.Cells.Item(156, 18).Value = "Hortaliza"

# Row 157 - "2a nueva(o)"
$ws.Cells.Item(157, 1).Value = 11
$ws.Cells.Item(157, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(157, 3).Value = "Bíobío"
$ws.Cells.Item(157, 4).Value = 44568
$ws.Cells.Item(157, 4).Style = $ws.Cells.Item(158, 4).Style
$ws.Cells.Item(157, 4).NumberFormat = $ws.Cells.Item(158, 4).NumberFormat
$ws.Cells.Item(157, 5).Value = 8
$ws.Cells.Item(157, 6).Value = 100112045
$ws.Cells.Item(157, 7).Value = "Zapallo"
$ws.Cells.Item(157, 8).Value = "Camote"
$ws.Cells.Item(157, 9).Value = "2a nueva(o)"
$ws.Cells.Item(157, 10).Value = 500
$ws.Cells.Item(157, 11).Value = 400
$ws.Cells.Item(157, 12).Value = 400
$ws.Cells.Item(157, 13).Value = 400
$ws.Cells.Item(157, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(157, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(157, 16).Value = 400
$ws.Cells.Item(157, 17).Value = 1
$ws.Cells.Item(157, 18).Value = "Hortaliza"
